$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 0.243056
$ws.Range("N2").Value = 0.729168
$ws.Range("O2").Value = 0.002199620488481675
$ws.Range("P2").Value = 0.002199620488481675
$ws.Range("Q2").Value = 48.01410141973867
$ws.Range("R2").Value = 432.126912777648
$ws.Range("S2").Value = 0.0007453862223379447
$ws.Range("T2").Value = 0.0007453862223379446
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("M3").Value = 70.95253000000001
$ws.Range("N3").Value = 212.85759
$ws.Range("O3").Value = 0.6421097964979703
$ws.Range("P3").Value = 0.6421097964979703
$ws.Range("Q3").Value = 14016.20190987694
$ws.Range("R3").Value = 126145.8171888925
$ws.Range("S3").Value = 0.2175919882743882
$ws.Range("T3").Value = 0.2175919882743881
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("M4").Value = 0.04794200000000001
$ws.Range("N4").Value = 0.143826
$ws.Range("O4").Value = 0.0004338679376719292
$ws.Range("P4").Value = 0.0004338679376719292
$ws.Range("Q4").Value = 9.470624260520667
$ws.Range("R4").Value = 85.235618344686
$ws.Range("S4").Value = 0.0001470249912420419
$ws.Range("T4").Value = 0.0001470249912420419
$ws.Range("G5").Value = 197.5433703333333
$ws.Range("H5").Value = 592.6301109999999
$ws.Range("I5").Value = 0.3388703761585983
$ws.Range("J5").Value = 0.3388703761585982
$ws.Range("M5").Value = 39.25553366666666
$ws.Range("N5").Value = 117.766601
$ws.Range("O5").Value = 0.3552567150758761
$ws.Range("P5").Value = 0.3552567150758761
$ws.Range("Q5").Value = 7754.670424746967
$ws.Range("R5").Value = 69792.0338227227
$ws.Range("S5").Value = 0.1203859766706301
$ws.Range("T5").Value = 0.1203859766706301
$ws.Range("G6").Value = 79.82725266666667
$ws.Range("H6").Value = 239.481758
$ws.Range("I6").Value = 0.1369374790620155
$ws.Range("J6").Value = 0.1369374790620154
$ws.Range("M6").Value = 0.243056
$ws.Range("N6").Value = 0.729168
$ws.Range("O6").Value = 0.002199620488481675
$ws.Range("P6").Value = 0.002199620488481675
$ws.Range("Q6").Value = 19.40249272414934
$ws.Range("R6").Value = 174.622434517344
$ws.Range("S6").Value = 0.0003012104845858396
$ws.Range("T6").Value = 0.0003012104845858395
$ws.Range("G7").Value = 79.82725266666667
$ws.Range("H7").Value = 239.481758
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("M7").Value = 70.95253000000001
$ws.Range("N7").Value = 212.85759
$ws.Range("O7").Value = 0.6421097964979703
$ws.Range("P7").Value = 0.6421097964979703
$ws.Range("Q7").Value = 5663.945539649248
$ws.Range("R7").Value = 50975.50985684322
$ws.Range("S7").Value = 0.08792889681345582
$ws.Range("T7").Value = 0.08792889681345581
$ws.Range("G8").Value = 79.82725266666667
$ws.Range("H8").Value = 239.481758
$ws.Range("I8").Value = 0.1369374790620155
$ws.Range("J8").Value = 0.1369374790620154
$ws.Range("M8").Value = 0.04794200000000001
$ws.Range("N8").Value = 0.143826
$ws.Range("O8").Value = 0.0004338679376719292
$ws.Range("P8").Value = 0.0004338679376719292
$ws.Range("Q8").Value = 3.827078147345334
$ws.Range("R8").Value = 34.443703326108
$ws.Range("S8").Value = 0.00005941278163062965
$ws.Range("T8").Value = 0.00005941278163062964
$ws.Range("G9").Value = 79.82725266666667
$ws.Range("H9").Value = 239.481758
$ws.Range("I9").Value = 0.1369374790620155
$ws.Range("J9").Value = 0.1369374790620154
$ws.Range("M9").Value = 39.25553366666666
$ws.Range("N9").Value = 117.766601
$ws.Range("O9").Value = 0.3552567150758761
$ws.Range("P9").Value = 0.3552567150758761
$ws.Range("Q9").Value = 3133.66140457384
$ws.Range("R9").Value = 28202.95264116456
$ws.Range("S9").Value = 0.04864795898234318
$ws.Range("T9").Value = 0.04864795898234318
$ws.Range("G10").Value = 148.824417
$ws.Range("H10").Value = 446.473251
$ws.Range("I10").Value = 0.2552967790580629
$ws.Range("J10").Value = 0.2552967790580629
$ws.Range("M10").Value = 0.243056
$ws.Range("N10").Value = 0.729168
$ws.Range("O10").Value = 0.002199620488481675
$ws.Range("P10").Value = 0.002199620488481675
$ws.Range("Q10").Value = 36.17266749835201
$ws.Range("R10").Value = 325.554007485168
$ws.Range("S10").Value = 0.0005615560258594944
$ws.Range("T10").Value = 0.0005615560258594944
$ws.Range("G11").Value = 148.824417
$ws.Range("H11").Value = 446.473251
$ws.Range("I11").Value = 0.2552967790580629
$ws.Range("J11").Value = 0.2552967790580629
$ws.Range("M11").Value = 70.95253000000001
$ws.Range("N11").Value = 212.85759
$ws.Range("O11").Value = 0.6421097964979703
$ws.Range("P11").Value = 0.6421097964979703
$ws.Range("Q11").Value = 10559.46891192501
$ws.Range("R11").Value = 95035.2202073251
$ws.Range("S11").Value = 0.16392856284756
$ws.Range("T11").Value = 0.16392856284756
$ws.Range("G12").Value = 148.824417
$ws.Range("H12").Value = 446.473251
$ws.Range("I12").Value = 0.2552967790580629
$ws.Range("J12").Value = 0.2552967790580629
$ws.Range("M12").Value = 0.04794200000000001
$ws.Range("N12").Value = 0.143826
$ws.Range("O12").Value = 0.0004338679376719292
$ws.Range("P12").Value = 0.0004338679376719292
$ws.Range("Q12").Value = 7.134940199814001
$ws.Range("R12").Value = 64.214461798326
$ws.Range("S12").Value = 0.0001107650870242079
$ws.Range("T12").Value = 0.0001107650870242079
$ws.Range("G13").Value = 148.824417
$ws.Range("H13").Value = 446.473251
$ws.Range("I13").Value = 0.2552967790580629
$ws.Range("J13").Value = 0.2552967790580629
$ws.Range("M13").Value = 39.25553366666666
$ws.Range("N13").Value = 117.766601
$ws.Range("O13").Value = 0.3552567150758761
$ws.Range("P13").Value = 0.3552567150758761
$ws.Range("Q13").Value = 5842.181911965539
$ws.Range("R13").Value = 52579.63720768985
$ws.Range("S13").Value = 0.09069589509761912
$ws.Range("T13").Value = 0.09069589509761913
$ws.Range("G14").Value = 35.426853
$ws.Range("H14").Value = 106.280559
$ws.Range("I14").Value = 0.06077202683121193
$ws.Range("J14").Value = 0.06077202683121192
$ws.Range("M14").Value = 0.243056
$ws.Range("N14").Value = 0.729168
$ws.Range("O14").Value = 0.002199620488481675
$ws.Range("P14").Value = 0.002199620488481675
$ws.Range("Q14").Value = 8.610709182768002
$ws.Range("R14").Value = 77.496382644912
$ws.Range("S14").Value = 0.0001336753953444918
$ws.Range("T14").Value = 0.0001336753953444918
$ws.Range("G15").Value = 35.426853
$ws.Range("H15").Value = 106.280559
$ws.Range("I15").Value = 0.06077202683121193
$ws.Range("J15").Value = 0.06077202683121192
$ws.Range("M15").Value = 70.95253000000001
$ws.Range("N15").Value = 212.85759
$ws.Range("O15").Value = 0.6421097964979703
$ws.Range("P15").Value = 0.6421097964979703
$ws.Range("Q15").Value = 2513.62485028809
$ws.Range("R15").Value = 22622.62365259281
$ws.Range("S15").Value = 0.03902231378135868
$ws.Range("T15").Value = 0.03902231378135868
$ws.Range("G16").Value = 35.426853
$ws.Range("H16").Value = 106.280559
$ws.Range("I16").Value = 0.06077202683121193
$ws.Range("J16").Value = 0.06077202683121192
$ws.Range("M16").Value = 0.04794200000000001
$ws.Range("N16").Value = 0.143826
$ws.Range("O16").Value = 0.0004338679376719292
$ws.Range("P16").Value = 0.0004338679376719292
$ws.Range("Q16").Value = 1.698434186526
$ws.Range("R16").Value = 15.285907678734
$ws.Range("S16").Value = 0.00002636703394940107
$ws.Range("T16").Value = 0.00002636703394940106
$ws.Range("G17").Value = 35.426853
$ws.Range("H17").Value = 106.280559
$ws.Range("I17").Value = 0.06077202683121193
$ws.Range("J17").Value = 0.06077202683121192
$ws.Range("M17").Value = 39.25553366666666
$ws.Range("N17").Value = 117.766601
$ws.Range("O17").Value = 0.3552567150758761
$ws.Range("P17").Value = 0.3552567150758761
$ws.Range("Q17").Value = 1390.700020645551
$ws.Range("R17").Value = 12516.30018580996
$ws.Range("S17").Value = 0.02158967062055935
$ws.Range("T17").Value = 0.02158967062055935
$ws.Range("G18").Value = 121.3248153333333
$ws.Range("H18").Value = 363.974446
$ws.Range("I18").Value = 0.2081233388901116
$ws.Range("J18").Value = 0.2081233388901115
$ws.Range("M18").Value = 0.243056
$ws.Range("N18").Value = 0.729168
$ws.Range("O18").Value = 0.002199620488481675
$ws.Range("P18").Value = 0.002199620488481675
$ws.Range("Q18").Value = 29.48872431565867
$ws.Range("R18").Value = 265.398518840928
$ws.Range("S18").Value = 0.0004577923603539043
$ws.Range("T18").Value = 0.0004577923603539042
$ws.Range("G19").Value = 121.3248153333333
$ws.Range("H19").Value = 363.974446
$ws.Range("I19").Value = 0.2081233388901116
$ws.Range("J19").Value = 0.2081233388901115
$ws.Range("M19").Value = 70.95253000000001
$ws.Range("N19").Value = 212.85759
$ws.Range("O19").Value = 0.6421097964979703
$ws.Range("P19").Value = 0.6421097964979703
$ws.Range("Q19").Value = 8608.302599682795
$ws.Range("R19").Value = 77474.72339714515
$ws.Range("S19").Value = 0.1336380347812076
$ws.Range("T19").Value = 0.1336380347812076
$ws.Range("G20").Value = 121.3248153333333
$ws.Range("H20").Value = 363.974446
$ws.Range("I20").Value = 0.2081233388901116
$ws.Range("J20").Value = 0.2081233388901115
$ws.Range("M20").Value = 0.04794200000000001
$ws.Range("N20").Value = 0.143826
$ws.Range("O20").Value = 0.0004338679376719292
$ws.Range("P20").Value = 0.0004338679376719292
$ws.Range("Q20").Value = 5.816554296710668
$ws.Range("R20").Value = 52.348988670396
$ws.Range("S20").Value = 0.00009029804382564873
$ws.Range("T20").Value = 0.00009029804382564871
$ws.Range("G21").Value = 121.3248153333333
$ws.Range("H21").Value = 363.974446
$ws.Range("I21").Value = 0.2081233388901116
$ws.Range("J21").Value = 0.2081233388901115
$ws.Range("M21").Value = 39.25553366666666
$ws.Range("N21").Value = 117.766601
$ws.Range("O21").Value = 0.3552567150758761
$ws.Range("P21").Value = 0.3552567150758761
$ws.Range("Q21").Value = 4762.670372919782
$ws.Range("R21").Value = 42864.03335627804
$ws.Range("S21").Value = 0.07393721370472436
$ws.Range("T21").Value = 0.07393721370472436
